$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44; this shifts the existing rows 44-143
# down to 45-144, matching the target diff (dimension A1:T143 -> A1:T144).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record.
$ws.Cells.Item(44, 1).Value = 10
$ws.Cells.Item(44, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(44, 3).Value = "La Araucanía"
$ws.Cells.Item(44, 4).Value = 44536
$ws.Cells.Item(44, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(44, 5).Value = 9
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100103
$ws.Cells.Item(44, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(44, 9).Value = 100103002
$ws.Cells.Item(44, 10).Value = "Ciruela"
$ws.Cells.Item(44, 11).Value = "Black Amber"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 95
$ws.Cells.Item(44, 14).Value = 18000
$ws.Cells.Item(44, 15).Value = 18000
$ws.Cells.Item(44, 16).Value = 18000
$ws.Cells.Item(44, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(44, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(44, 19).Value = 1000
$ws.Cells.Item(44, 20).Value = 18
